$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row labels from the "_old"/"_new" suffix convention
#    to the "_FV2210"/"_FV2304" (format-version) suffix convention.
$oldHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$newHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i]
}
# column 11 ("K") keeps its "diff" label
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
}

# 2) Freeze the header row.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the data range into an Excel table ("ListObject") so the header
#    row can be filtered/sorted.
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $dataRange, $false, 1)
$tbl.Name = "Table1"
